# Apply updated dSF (column F) values as per repull/push of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -11
    3  = -4
    4  = -2
    6  = -2
    8  = -1
    10 = -6
    11 = -8
    13 = -1
    14 = -3
    18 = 5
    19 = -3
    21 = -3
    23 = -1
    24 = -5
    28 = 1
    30 = -4
    34 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
